# Generate Report for Handback
# Marks the zh-cn / de-de / Overview "Status" as handed back (in sync with
# en-US), and fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale detail sheets now
# that the handback xliffs have been generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$mdFile     = "fa70e7e0-2a4a-4545-bca5-1f19d4af4414.md"

# ---- Overview sheet: zh-cn / de-de status columns ----
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns("E").ColumnWidth = 29.9777047293527
$overview.Columns("F").ColumnWidth = 29.9777047293527

# ---- zh-cn detail sheet ----
$zhcn.Range("C2").Value = $statusText

# Latest Target File (I2) now points at the source .md, same as A2, styled
# (and linked) like a hyperlink.
$zhcn.Range("I2").Value = $zhcn.Range("A2").Value2
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/aa607e9ca084e6c63cc6fbec476007fffdaec3fb/e2e/$mdFile", [Type]::Missing, $mdFile, $mdFile) | Out-Null

# Latest Handback File (J2) mirrors the Latest Handoff File (G2) - the
# handback xliff uses the same name as the handoff xliff.
$zhcn.Range("J2").Value = $zhcn.Range("G2").Value2

# Latest Handback DateTime (K2) - generated just now.
$zhcn.Range("K2").Value = "2016-08-13 11:12:52"

$zhcn.Columns("C").ColumnWidth = 29.9777047293527
$zhcn.Columns("I").ColumnWidth = 40
$zhcn.Columns("J").ColumnWidth = 40

# ---- de-de detail sheet ----
$dede.Range("C2").Value = $statusText

$dede.Range("I2").Value = $dede.Range("A2").Value2
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/aa607e9ca084e6c63cc6fbec476007fffdaec3fb/e2e/$mdFile", [Type]::Missing, $mdFile, $mdFile) | Out-Null

$dede.Range("J2").Value = $dede.Range("G2").Value2

$dede.Range("K2").Value = "2016-08-13 11:13:04"

$dede.Columns("C").ColumnWidth = 29.9777047293527
$dede.Columns("I").ColumnWidth = 40
$dede.Columns("J").ColumnWidth = 40
